$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "64.664.34"
$ws.Cells.Item(2, 5).Value = "  -0.71%  "
$ws.Cells.Item(3, 4).Value = "3.423.46"
$ws.Cells.Item(3, 5).Value = "  -1.65%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "572.81"
$ws.Cells.Item(5, 5).Value = "  -1.04%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "156.93"
$ws.Cells.Item(6, 5).Value = "  -3.34%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.620"
$ws.Cells.Item(7, 5).Value = "  +6.37%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.00"
$ws.Cells.Item(8, 5).Value = "  +0.10%  "
$ws.Cells.Item(9, 4).Value = "3.430.06"
$ws.Cells.Item(9, 5).Value = "  -1.48%  "
$ws.Cells.Item(10, 5).Value = "  -2.40%  "
$ws.Cells.Item(11, 5).Value = "  -2.68%  "
$ws.Cells.Item(12, 5).Value = "  +0.07%  "
$ws.Cells.Item(13, 4).Value = "4.014.80"
$ws.Cells.Item(13, 5).Value = "  -1.32%  "
$ws.Cells.Item(14, 5).Value = "  +0.22%  "
$ws.Cells.Item(15, 5).Value = "  -3.69%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "27.98"
$ws.Cells.Item(16, 5).Value = "  -3.23%  "
$ws.Cells.Item(17, 4).Value = "64.685.21"
$ws.Cells.Item(17, 5).Value = "  -0.59%  "
$ws.Cells.Item(18, 4).Value = "3.369.88"
$ws.Cells.Item(18, 5).Value = "  -2.24%  "
$ws.Cells.Item(19, 5).Value = "  -0.27%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.02"
$ws.Cells.Item(20, 5).Value = "  -2.46%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "377.00"
$ws.Cells.Item(21, 5).Value = "  -3.69%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "8.05"
$ws.Cells.Item(22, 5).Value = "  -2.40%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.553"
$ws.Cells.Item(23, 5).Value = "  +0.52%  "
$ws.Cells.Item(24, 5).Value = "  -0.62%  "
$ws.Cells.Item(25, 5).Value = "  -0.80%  "
$ws.Cells.Item(26, 5).Value = "  -4.70%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.24"
$ws.Cells.Item(27, 5).Value = "  +7.03%  "
$ws.Cells.Item(28, 5).Value = "  -2.33%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.997"
$ws.Cells.Item(29, 5).Value = "  -0.26%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.49"
$ws.Cells.Item(30, 5).Value = "  +2.71%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.22"
$ws.Cells.Item(31, 5).Value = "  +0.28%  "
$ws.Cells.Item(32, 5).Value = "  -0.93%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "23.13"
$ws.Cells.Item(33, 5).Value = "  -2.55%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "7.25"
$ws.Cells.Item(34, 5).Value = "  +2.02%  "
$ws.Cells.Item(35, 5).Value = "  +6.81%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "160.28"
$ws.Cells.Item(36, 5).Value = "  -0.91%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.90"
$ws.Cells.Item(37, 5).Value = "  -0.48%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "6.98"
$ws.Cells.Item(38, 5).Value = "  +7.09%  "
$ws.Cells.Item(39, 5).Value = "  -1.11%  "
$ws.Cells.Item(40, 4).Value = "2.884.01"
$ws.Cells.Item(40, 5).Value = "  -3.97%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "26.76"
$ws.Cells.Item(41, 5).Value = "  -2.76%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "4.62"
$ws.Cells.Item(42, 5).Value = "  +0.75%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "26.67"
$ws.Cells.Item(43, 5).Value = "  +9.69%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "42.77"
$ws.Cells.Item(44, 5).Value = "  -0.34%  "
$ws.Cells.Item(45, 5).Value = "  -1.04%  "
$ws.Cells.Item(46, 5).Value = "  -1.14%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "321.24"
$ws.Cells.Item(47, 5).Value = "  +5.24%  "
$ws.Cells.Item(48, 5).Value = "  -1.97%  "
$ws.Cells.Item(49, 5).Value = "  +1.26%  "
$ws.Cells.Item(50, 5).Value = "  +0.61%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.16"
$ws.Cells.Item(51, 5).Value = "  -1.88%  "
